$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these cells to Text format first so numeric-looking strings
# (e.g. "238.40", "1.001") are stored as text, matching the source data,
# instead of being auto-coerced to numbers (which would also drop trailing zeros).
$textCells = @("D5", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = '25.893.66'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '1.740.79'
$ws.Range("E3").Value = '  +0.48%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '238.40'
$ws.Range("E5").Value = '  +4.28%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '0.5173'
$ws.Range("E7").Value = '  -1.08%  '
$ws.Range("D8").Value = '0.2749'
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '0.06156'
$ws.Range("E9").Value = '  +0.59%  '
$ws.Range("D10").Value = '1.741.99'
$ws.Range("E10").Value = '  +0.61%  '
$ws.Range("D11").Value = '0.07170'
$ws.Range("E11").Value = '  +1.51%  '
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '14.99'
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '0.6445'
$ws.Range("E13").Value = '  +1.55%  '
$ws.Range("D14").Value = '4.602'
$ws.Range("E14").Value = '  +1.80%  '
$ws.Range("D15").Value = '77.45'
$ws.Range("E15").Value = '  +1.23%  '
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("D18").Value = '25.911.96'
$ws.Range("E18").Value = '  +0.37%  '
$ws.Range("E19").Value = '  +2.33%  '
$ws.Range("D20").Value = '0.000006776'
$ws.Range("E20").Value = '  +2.24%  '
$ws.Range("D21").Value = '1.966.98'
$ws.Range("E21").Value = '  +0.35%  '
$ws.Range("D22").Value = '4.278'
$ws.Range("D23").Value = '8.656'
$ws.Range("E23").Value = '  -1.27%  '
$ws.Range("D24").Value = '5.258'
$ws.Range("E24").Value = '  +1.88%  '
$ws.Range("D25").Value = '139.08'
$ws.Range("E25").Value = '  -0.69%  '
$ws.Range("D26").Value = '1.511'
$ws.Range("E26").Value = '  +0.44%  '
$ws.Range("D27").Value = '15.13'
$ws.Range("E27").Value = '  +0.88%  '
$ws.Range("D28").Value = '1.761'
$ws.Range("E28").Value = '  -0.73%  '
$ws.Range("D29").Value = '105.86'
$ws.Range("E29").Value = '  +3.78%  '
$ws.Range("D30").Value = '3.942'
$ws.Range("E30").Value = '  +6.36%  '
$ws.Range("D31").Value = '0.08292'
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("D32").Value = '3.685'
$ws.Range("E32").Value = '  +5.33%  '
$ws.Range("D33").Value = '0.04595'
$ws.Range("E33").Value = '  +3.33%  '
$ws.Range("D34").Value = '2.642'
$ws.Range("E34").Value = '  +1.19%  '
$ws.Range("D35").Value = '0.9874'
$ws.Range("E35").Value = '  +2.05%  '
$ws.Range("D36").Value = '0.6179'
$ws.Range("E36").Value = '  +0.39%  '
$ws.Range("D37").Value = '2.683'
$ws.Range("E37").Value = '  +0.59%  '
$ws.Range("E38").Value = '  +3.05%  '
$ws.Range("D39").Value = '1.927'
$ws.Range("E39").Value = '  +1.50%  '
$ws.Range("D40").Value = '0.9999'
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").Value = '97.83'
$ws.Range("D42").Value = '0.3833'
$ws.Range("E42").Value = '  +0.64%  '
$ws.Range("D43").Value = '0.7413'
$ws.Range("E43").Value = '  +3.02%  '
$ws.Range("D44").Value = '4.987'
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("D45").Value = '0.1129'
$ws.Range("E45").Value = '  +1.11%  '
$ws.Range("D46").Value = '6.201'
$ws.Range("E46").Value = '  +1.03%  '
$ws.Range("D47").Value = '0.05258'
$ws.Range("E47").Value = '  -1.33%  '
$ws.Range("D48").Value = '55.00'
$ws.Range("E48").Value = '  +3.50%  '
$ws.Range("E49").Value = '  +1.94%  '
$ws.Range("D50").Value = '7.628'
$ws.Range("E50").Value = '  +0.45%  '
$ws.Range("D51").Value = '0.3404'
$ws.Range("E51").Value = '  +1.26%  '

Write-Host "Applied 96 cell updates"
